$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 161.375
$ws.Range("I5").Value = 138
$ws.Range("K5").Value = 138
$ws.Range("M5").Value = -23

# Row 64
$ws.Range("H64").Value = 3442.623
$ws.Range("J64").Value = 4733.3335
$ws.Range("L64").Value = 4733.3335
$ws.Range("N64").Value = -5229.3335

# Row 67
$ws.Range("H67").Value = 3442.623
$ws.Range("J67").Value = 4733.3335
$ws.Range("L67").Value = 4733.3335
$ws.Range("N67").Value = -6449.3335

# Row 92
$ws.Range("H92").Value = 1021.35
$ws.Range("I92").Value = 245.4375
$ws.Range("K92").Value = 245.4375
$ws.Range("M92").Value = 1002.5625

# Row 132
$ws.Range("H132").Value = 4397.091
$ws.Range("I132").Value = 4263.1113
$ws.Range("K132").Value = 12789.3339
$ws.Range("M132").Value = -10259.3339

# Row 138
$ws.Range("H138").Value = 3671.6785
$ws.Range("J138").Value = 3568.84
$ws.Range("L138").Value = 10706.52
$ws.Range("N138").Value = -20986.52


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7702.963
$ws.Range("I32").Value = 3173.1272
$ws.Range("K32").Value = 3173.1272
$ws.Range("M32").Value = -2886.1272

# Row 61
$ws.Range("H61").Value = 3230125.5
$ws.Range("I61").Value = 4617.615
$ws.Range("K61").Value = 4617.615
$ws.Range("M61").Value = -4405.615

# Row 74
$ws.Range("H74").Value = 1639858.9
$ws.Range("I74").Value = 2418012.5
$ws.Range("J74").Value = 12810.272
$ws.Range("K74").Value = 2418012.5
$ws.Range("L74").Value = 12810.272
$ws.Range("M74").Value = -2417138.5
$ws.Range("N74").Value = -14558.272

# Row 77
$ws.Range("H77").Value = 1639858.9
$ws.Range("I77").Value = 2418012.5
$ws.Range("J77").Value = 12810.272
$ws.Range("K77").Value = 12090062.5
$ws.Range("L77").Value = 64051.36
$ws.Range("M77").Value = -12085694.5
$ws.Range("N77").Value = -72787.36

# Row 102
$ws.Range("H102").Value = 2013.4286
$ws.Range("I102").Value = 2299.0833
$ws.Range("K102").Value = 2299.0833
$ws.Range("M102").Value = -677.0832999999998

# Row 110
$ws.Range("H110").Value = 1648.4
$ws.Range("I110").Value = 1123.75
$ws.Range("K110").Value = 1123.75
$ws.Range("M110").Value = 921.25

# Row 132
$ws.Range("H132").Value = 809225.5600000001
$ws.Range("I132").Value = 1193485.5
$ws.Range("J132").Value = 2279.8
$ws.Range("K132").Value = 3580456.5
$ws.Range("L132").Value = 6839.400000000001
$ws.Range("M132").Value = -3577926.5
$ws.Range("N132").Value = -11899.4

# Row 136
$ws.Range("H136").Value = 3230125.5
$ws.Range("I136").Value = 4617.615
$ws.Range("K136").Value = 13852.845
$ws.Range("M136").Value = -11302.845

# Row 139
$ws.Range("H139").Value = 119182.86
$ws.Range("J139").Value = 119182.86
$ws.Range("L139").Value = 119182.86
$ws.Range("N139").Value = -129462.86


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 36679.566
$ws.Range("I20").Value = 47350.78
$ws.Range("K20").Value = 47350.78
$ws.Range("M20").Value = -47103.78

# Row 36
$ws.Range("H36").Value = 899
$ws.Range("I36").Value = 899
$ws.Range("K36").Value = 899
$ws.Range("M36").Value = -365

# Row 86
$ws.Range("H86").Value = 1804.8
$ws.Range("I86").Value = 1804.8
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1804.8
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -681.8
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 1804.8
$ws.Range("I89").Value = 1804.8
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9024
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3408
$ws.Range("N89").ClearContents()

# Row 134
$ws.Range("H134").Value = 7948293.5
$ws.Range("I134").Value = 8250
$ws.Range("J134").Value = 23828380
$ws.Range("K134").Value = 24750
$ws.Range("L134").Value = 71485140
$ws.Range("M134").Value = -22215
$ws.Range("N134").Value = -71490210


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 17749.5

# Row 50
$ws.Range("H50").Value = 18350

# Row 51
$ws.Range("H51").Value = 24350

# Row 58
$ws.Range("H58").Value = 3797981.5
$ws.Range("I58").Value = 7252.2
$ws.Range("J58").Value = 4912902
$ws.Range("K58").Value = 7252.2
$ws.Range("L58").Value = 4912902
$ws.Range("M58").Value = -7049.2
$ws.Range("N58").Value = -4913308

# Row 60
$ws.Range("H60").Value = 19402.75

# Row 61
$ws.Range("H61").Value = 24350

# Row 69
$ws.Range("H69").Value = 10061
$ws.Range("I69").Value = 4355.4287
$ws.Range("J69").Value = 50000
$ws.Range("K69").Value = 4355.4287
$ws.Range("L69").Value = 50000
$ws.Range("M69").Value = -3606.4287
$ws.Range("N69").Value = -51498

# Row 72
$ws.Range("H72").Value = 10061
$ws.Range("I72").Value = 4355.4287
$ws.Range("J72").Value = 50000
$ws.Range("K72").Value = 13066.2861
$ws.Range("L72").Value = 150000
$ws.Range("M72").Value = -9322.286100000001
$ws.Range("N72").Value = -157488

# Row 99
$ws.Range("H99").Value = 35140.555
$ws.Range("I99").Value = 33569.715
$ws.Range("K99").Value = 33569.715
$ws.Range("M99").Value = -32071.715

# Row 105
$ws.Range("H105").Value = 10274.75
$ws.Range("I105").Value = 12512.111
$ws.Range("K105").Value = 12512.111
$ws.Range("M105").Value = -10765.111

# Row 126
$ws.Range("H126").Value = 35140.555
$ws.Range("I126").Value = 33569.715
$ws.Range("K126").Value = 100709.145
$ws.Range("M126").Value = -98239.14499999999

# Row 132
$ws.Range("H132").Value = 3360.7778
$ws.Range("I132").Value = 3264.4707
$ws.Range("K132").Value = 9793.4121
$ws.Range("M132").Value = -7263.4121

# Row 134
$ws.Range("H134").Value = 2499.5
$ws.Range("I134").Value = 3999
$ws.Range("K134").Value = 11997
$ws.Range("M134").Value = -9462

# Row 136
$ws.Range("H136").Value = 3797981.5
$ws.Range("I136").Value = 7252.2
$ws.Range("J136").Value = 4912902
$ws.Range("K136").Value = 21756.6
$ws.Range("L136").Value = 14738706
$ws.Range("M136").Value = -19206.6
$ws.Range("N136").Value = -14743806


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 950
$ws.Range("I7").Value = 950
$ws.Range("K7").Value = 2850
$ws.Range("M7").Value = -2738

# Row 12
$ws.Range("H12").Value = 802.5714
$ws.Range("I12").Value = 1292.5
$ws.Range("J12").Value = 149.33333
$ws.Range("K12").Value = 3877.5
$ws.Range("L12").Value = 447.99999
$ws.Range("M12").Value = -3704.5
$ws.Range("N12").Value = -793.99999

# Row 34
$ws.Range("H34").Value = 339.75
$ws.Range("I34").Value = 29.5
$ws.Range("K34").Value = 88.5
$ws.Range("M34").Value = -4.5

# Row 92
$ws.Range("H92").Value = 355.58823
$ws.Range("I92").Value = 364.81818
$ws.Range("J92").Value = 338.66666
$ws.Range("K92").Value = 1094.45454
$ws.Range("L92").Value = 1015.99998
$ws.Range("M92").Value = 153.54546
$ws.Range("N92").Value = -3511.99998

# Row 96
$ws.Range("H96").Value = 14337.333
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -49118

# Row 102
$ws.Range("H102").Value = 4283.3335
$ws.Range("I102").Value = 4283.3335
$ws.Range("K102").Value = 12850.0005
$ws.Range("M102").Value = -10416.0005

# Row 113
$ws.Range("H113").Value = 1418.7693
$ws.Range("I113").Value = 745.25
$ws.Range("J113").Value = 1718.1111
$ws.Range("K113").Value = 2235.75
$ws.Range("L113").Value = 5154.3333
$ws.Range("M113").Value = -65.75
$ws.Range("N113").Value = -9494.3333

# Row 138
$ws.Range("H138").Value = 23262.824
$ws.Range("I138").Value = 23262.824
$ws.Range("K138").Value = 69788.47200000001
$ws.Range("M138").Value = -64648.47200000001


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 7000
$ws.Range("J5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("N5").Value = -7224

# Row 26
$ws.Range("H26").Value = 14999.5
$ws.Range("J26").Value = 15000
$ws.Range("L26").Value = 15000
$ws.Range("N26").Value = -15560

# Row 50
$ws.Range("H50").Value = 14999.5
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -15996

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 122
$ws.Range("H122").Value = 47615.75
$ws.Range("I122").Value = 60765.61
$ws.Range("K122").Value = 182296.83
$ws.Range("M122").Value = -179846.83


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1108.55
$ws.Range("I16").Value = 330.57144
$ws.Range("J16").Value = 6554.4
$ws.Range("K16").Value = 330.57144
$ws.Range("L16").Value = 6554.4
$ws.Range("M16").Value = -160.57144
$ws.Range("N16").Value = -6894.4

# Row 74
$ws.Range("H74").Value = 49808
$ws.Range("I74").Value = 49808
$ws.Range("K74").Value = 49808
$ws.Range("M74").Value = -48810

# Row 77
$ws.Range("H77").Value = 49808
$ws.Range("I77").Value = 49808
$ws.Range("K77").Value = 149424
$ws.Range("M77").Value = -144432

# Row 122
$ws.Range("H122").Value = 4478.6
$ws.Range("I122").Value = 2578.6
$ws.Range("K122").Value = 7735.799999999999
$ws.Range("M122").Value = -5285.799999999999

# Row 132
$ws.Range("H132").Value = 5314272
$ws.Range("I132").Value = 8990718
$ws.Range("J132").Value = 3850
$ws.Range("K132").Value = 26972154
$ws.Range("L132").Value = 11550
$ws.Range("M132").Value = -26969624
$ws.Range("N132").Value = -16610

# Row 134
$ws.Range("H134").Value = 55000
$ws.Range("J134").Value = 55000
$ws.Range("L134").Value = 55000
$ws.Range("N134").Value = -65140

# Row 136
$ws.Range("H136").Value = 6103189.5
$ws.Range("I136").Value = 9618381
$ws.Range("K136").Value = 28855143
$ws.Range("M136").Value = -28852593

